$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the date-like text columns (Y, AA) as Text so the values are
# stored as literal strings "2023-07-27" rather than being auto-parsed into
# date serial numbers. ClearFormats afterwards removes the temporary style
# so the cells end up with no explicit style, matching the source data.
# (NumberFormat/ClearFormats must be applied per single-area range - a
# multi-area "A,B" reference only affects the first area in this engine.)
$ws.Range("Y8:Y15").NumberFormat = "@"
$ws.Range("AA8:AA15").NumberFormat = "@"

# Row 8
$ws.Range("A8").Value = 112243565
$ws.Range("B8").Value = 90666
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = 'Dropptaggsvamp'
$ws.Range("G8").Value = 'Hydnellum ferrugineum'
$ws.Range("H8").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("I8").Value = ""
$ws.Range("P8").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q8").Value = 404459
$ws.Range("R8").Value = 6706753
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 'Värmland'
$ws.Range("U8").Value = 'Torsby'
$ws.Range("V8").Value = 'Värmland'
$ws.Range("W8").Value = 'Norra Ny'
$ws.Range("Y8").Value = '2023-07-27'
$ws.Range("AA8").Value = '2023-07-27'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Value = ""
$ws.Range("AW8").Value = 'Dick Östberg'
$ws.Range("AX8").Value = 'Dick Östberg'
$ws.Range("AY8").Value = ""

# Row 9
$ws.Range("A9").Value = 112243563
$ws.Range("B9").Value = 89425
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 5442
$ws.Range("F9").Value = 'Tallticka'
$ws.Range("G9").Value = 'Porodaedalea pini'
$ws.Range("H9").Value = '(Brot.) Murrill'
$ws.Range("I9").Value = ""
$ws.Range("P9").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q9").Value = 404744
$ws.Range("R9").Value = 6707084
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Värmland'
$ws.Range("U9").Value = 'Torsby'
$ws.Range("V9").Value = 'Värmland'
$ws.Range("W9").Value = 'Norra Ny'
$ws.Range("Y9").Value = '2023-07-27'
$ws.Range("AA9").Value = '2023-07-27'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AT9").Value = ""
$ws.Range("AW9").Value = 'Dick Östberg'
$ws.Range("AX9").Value = 'Dick Östberg'
$ws.Range("AY9").Value = ""

# Row 10
$ws.Range("A10").Value = 112243594
$ws.Range("B10").Value = 77515
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = 'Garnlav'
$ws.Range("G10").Value = 'Alectoria sarmentosa'
$ws.Range("H10").Value = '(Ach.) Ach.'
$ws.Range("I10").Value = ""
$ws.Range("P10").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q10").Value = 404742
$ws.Range("R10").Value = 6706992
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Värmland'
$ws.Range("U10").Value = 'Torsby'
$ws.Range("V10").Value = 'Värmland'
$ws.Range("W10").Value = 'Norra Ny'
$ws.Range("Y10").Value = '2023-07-27'
$ws.Range("AA10").Value = '2023-07-27'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").Value = ""
$ws.Range("AW10").Value = 'Dick Östberg'
$ws.Range("AX10").Value = 'Dick Östberg'
$ws.Range("AY10").Value = ""

# Row 11
$ws.Range("A11").Value = 112243589
$ws.Range("B11").Value = 77515
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = 'Garnlav'
$ws.Range("G11").Value = 'Alectoria sarmentosa'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("I11").Value = ""
$ws.Range("P11").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q11").Value = 404762
$ws.Range("R11").Value = 6707097
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Värmland'
$ws.Range("U11").Value = 'Torsby'
$ws.Range("V11").Value = 'Värmland'
$ws.Range("W11").Value = 'Norra Ny'
$ws.Range("Y11").Value = '2023-07-27'
$ws.Range("AA11").Value = '2023-07-27'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AT11").Value = ""
$ws.Range("AW11").Value = 'Dick Östberg'
$ws.Range("AX11").Value = 'Dick Östberg'
$ws.Range("AY11").Value = ""

# Row 12
$ws.Range("A12").Value = 112243573
$ws.Range("B12").Value = 77515
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("I12").Value = ""
$ws.Range("P12").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q12").Value = 404477
$ws.Range("R12").Value = 6706766
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Värmland'
$ws.Range("U12").Value = 'Torsby'
$ws.Range("V12").Value = 'Värmland'
$ws.Range("W12").Value = 'Norra Ny'
$ws.Range("Y12").Value = '2023-07-27'
$ws.Range("AA12").Value = '2023-07-27'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AT12").Value = ""
$ws.Range("AW12").Value = 'Dick Östberg'
$ws.Range("AX12").Value = 'Dick Östberg'
$ws.Range("AY12").Value = ""

# Row 13
$ws.Range("A13").Value = 112243569
$ws.Range("B13").Value = 77515
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = 'Garnlav'
$ws.Range("G13").Value = 'Alectoria sarmentosa'
$ws.Range("H13").Value = '(Ach.) Ach.'
$ws.Range("I13").Value = ""
$ws.Range("P13").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q13").Value = 404751
$ws.Range("R13").Value = 6707073
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Värmland'
$ws.Range("U13").Value = 'Torsby'
$ws.Range("V13").Value = 'Värmland'
$ws.Range("W13").Value = 'Norra Ny'
$ws.Range("Y13").Value = '2023-07-27'
$ws.Range("AA13").Value = '2023-07-27'
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AT13").Value = ""
$ws.Range("AW13").Value = 'Dick Östberg'
$ws.Range("AX13").Value = 'Dick Östberg'
$ws.Range("AY13").Value = ""

# Row 14
$ws.Range("A14").Value = 112243600
$ws.Range("B14").Value = 77515
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = 'Garnlav'
$ws.Range("G14").Value = 'Alectoria sarmentosa'
$ws.Range("H14").Value = '(Ach.) Ach.'
$ws.Range("I14").Value = ""
$ws.Range("P14").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q14").Value = 404725
$ws.Range("R14").Value = 6707036
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 'Värmland'
$ws.Range("U14").Value = 'Torsby'
$ws.Range("V14").Value = 'Värmland'
$ws.Range("W14").Value = 'Norra Ny'
$ws.Range("Y14").Value = '2023-07-27'
$ws.Range("AA14").Value = '2023-07-27'
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AT14").Value = ""
$ws.Range("AW14").Value = 'Dick Östberg'
$ws.Range("AX14").Value = 'Dick Östberg'
$ws.Range("AY14").Value = ""

# Row 15
$ws.Range("A15").Value = 112243588
$ws.Range("B15").Value = 77515
$ws.Range("C15").Value = 'Ovaliderad'
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = 'Garnlav'
$ws.Range("G15").Value = 'Alectoria sarmentosa'
$ws.Range("H15").Value = '(Ach.) Ach.'
$ws.Range("I15").Value = ""
$ws.Range("P15").Value = 'Väst Värsjön, Vrm'
$ws.Range("Q15").Value = 404452
$ws.Range("R15").Value = 6706739
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = 'Värmland'
$ws.Range("U15").Value = 'Torsby'
$ws.Range("V15").Value = 'Värmland'
$ws.Range("W15").Value = 'Norra Ny'
$ws.Range("Y15").Value = '2023-07-27'
$ws.Range("AA15").Value = '2023-07-27'
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AT15").Value = ""
$ws.Range("AW15").Value = 'Dick Östberg'
$ws.Range("AX15").Value = 'Dick Östberg'
$ws.Range("AY15").Value = ""

# Remove the temporary Text style from the date-like columns now that the
# values have been safely written as literal strings.
$ws.Range("Y8:Y15,AA8:AA15").ClearFormats()
